# Updated symbol list on Tue Dec 20 13:42:26 UTC 2022 with GitHub Actions
# Refreshes the cryptocurrency price snapshot (column D) and a couple of
# mislabelled coin-name cells (column E) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free approach: column D values are stored as *text* (not numbers)
# in this workbook, so each new price is entered with a leading apostrophe
# (forces Excel to keep it as text instead of auto-converting to a Number)
# and then the cell style is reset back to "Normal" so the quote-prefix
# formatting introduced by the apostrophe entry doesn't linger.

function Set-TextValue($range, $text) {
    $range.Formula = "'" + $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2")  "249.19"
Set-TextValue $ws.Range("D3")  "22.94"
Set-TextValue $ws.Range("D4")  "5.390"
Set-TextValue $ws.Range("D5")  "0.05617"
Set-TextValue $ws.Range("D6")  "3.441"
Set-TextValue $ws.Range("D7")  "6.358"
Set-TextValue $ws.Range("D8")  "0.8179"
Set-TextValue $ws.Range("D9")  "0.9174"
Set-TextValue $ws.Range("D10") "0.1434"
Set-TextValue $ws.Range("D11") "0.07506"
Set-TextValue $ws.Range("D12") "0.03195"
Set-TextValue $ws.Range("D13") "0.03095"
Set-TextValue $ws.Range("D14") "0.09325"
Set-TextValue $ws.Range("D15") "3.556"
Set-TextValue $ws.Range("D16") "0.001633"
Set-TextValue $ws.Range("D17") "0.04717"
Set-TextValue $ws.Range("D18") "0.0005764"

$ws.Range("E18").Value = "17OneONEWorstin24h"

Set-TextValue $ws.Range("D19") "0.006406"
Set-TextValue $ws.Range("D20") "0.005056"
Set-TextValue $ws.Range("D22") "0.0001500"
Set-TextValue $ws.Range("D23") "3.727"
Set-TextValue $ws.Range("D24") "2.163"
Set-TextValue $ws.Range("D25") "0.3296"

$ws.Range("E27").Value = "26AAXTokenAAB"

Set-TextValue $ws.Range("D40") "0.04002"
Set-TextValue $ws.Range("D41") "0.006880"
Set-TextValue $ws.Range("D42") "0.1068"
Set-TextValue $ws.Range("D43") "0.002760"
Set-TextValue $ws.Range("D44") "0.008578"
Set-TextValue $ws.Range("D45") "0.00005576"
Set-TextValue $ws.Range("D46") "0.00000000751"
Set-TextValue $ws.Range("D49") "0.2089"
Set-TextValue $ws.Range("D51") "0.01011"
